$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price values in column D are written as text,
# matching the inlineStr/shared-string cell type used in the workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.169.43'
$ws.Range("E2").Value = '  +5.53%  '
$ws.Range("D3").Value = '1.779.11'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '243.43'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("D7").Value = '0.4899'
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("D8").Value = '0.2661'
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").Value = '1.775.70'
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("E11").Value = '  +4.04%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = '0.6261'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").Value = '4.627'
$ws.Range("E14").Value = '  +3.00%  '
$ws.Range("D15").Value = '79.68'
$ws.Range("E15").Value = '  +3.21%  '
$ws.Range("D16").Value = '28.152.13'
$ws.Range("E16").Value = '  +6.24%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '0.000007223'
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").Value = '12.03'
$ws.Range("E20").Value = '  +5.67%  '
$ws.Range("D21").Value = '2.006.73'
$ws.Range("E21").Value = '  +2.78%  '
$ws.Range("D22").Value = '4.558'
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("D23").Value = '8.714'
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").Value = '5.212'
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("E26").Value = '  +2.95%  '
$ws.Range("D27").Value = '1.856'
$ws.Range("E27").Value = '  +4.84%  '
$ws.Range("D28").Value = '109.07'
$ws.Range("E28").Value = '  +2.74%  '
$ws.Range("D29").Value = '1.378'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '4.169'
$ws.Range("E30").Value = '  +6.32%  '
$ws.Range("D31").Value = '0.08230'
$ws.Range("E31").Value = '  +3.26%  '
$ws.Range("D32").Value = '3.752'
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("E33").Value = '  +9.22%  '
$ws.Range("E34").Value = '  +6.66%  '
$ws.Range("D35").Value = '2.614'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = '0.6489'
$ws.Range("D37").Value = '0.9473'
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("D38").Value = '2.589'
$ws.Range("E38").Value = '  +7.50%  '
$ws.Range("D39").Value = '2.039'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").Value = '5.882'
$ws.Range("E40").Value = '  +5.03%  '
$ws.Range("D41").Value = '0.01544'
$ws.Range("E41").Value = '  +2.17%  '
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").Value = '99.49'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '0.3966'
$ws.Range("E44").Value = '  +3.04%  '
$ws.Range("D45").Value = '7.145'
$ws.Range("E45").Value = '  +3.91%  '
$ws.Range("E46").Value = '  +4.10%  '
$ws.Range("D47").Value = '0.05430'
$ws.Range("E47").Value = '  +0.94%  '
$ws.Range("D48").Value = '7.993'
$ws.Range("E48").Value = '  +2.12%  '
$ws.Range("D49").Value = '1.295'
$ws.Range("E49").Value = '  +5.43%  '
$ws.Range("D50").Value = '30.58'
$ws.Range("E50").Value = '  +1.23%  '
$ws.Range("E51").Value = '  +2.14%  '

# Restore default styling on column D (removes the temporary text format
# without leaving a residual style reference on the cells).
$priceRange.ClearFormats()
